$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 210, shifting existing rows 210:316 down to 211:317.
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new record's data.
# (Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T mirror the previous row 210 values;
#  D,M,N,O,P,S carry the new figures per the commit.)
$ws.Cells.Item(210, 1).Value = 10
$ws.Cells.Item(210, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(210, 3).Value = "La Araucanía"
$ws.Cells.Item(210, 4).Value = 44839
$ws.Cells.Item(210, 5).Value = 9
$ws.Cells.Item(210, 6).Value = "Fruta"
$ws.Cells.Item(210, 7).Value = 100102
$ws.Cells.Item(210, 8).Value = "Cítricos"
$ws.Cells.Item(210, 9).Value = 100102006
$ws.Cells.Item(210, 10).Value = "Pomelo"
$ws.Cells.Item(210, 11).Value = "Start Ruby"
$ws.Cells.Item(210, 12).Value = "Primera"
$ws.Cells.Item(210, 13).Value = 95
$ws.Cells.Item(210, 14).Value = 14000
$ws.Cells.Item(210, 15).Value = 14000
$ws.Cells.Item(210, 16).Value = 14000
$ws.Cells.Item(210, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(210, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(210, 19).Value = 933
$ws.Cells.Item(210, 20).Value = 15
